$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting existing rows 103:176 down to 104:177.
$ws.Rows("103:103").Insert()

# Populate the newly inserted row 103 with a new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T repeat the same values used throughout the sheet
# for this market/product; D,L,M,N,O,P,S hold the new week's data.
$ws.Range("A103").Value = 1
$ws.Range("B103").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C103").Value = "Arica y Parinacota"
$ws.Range("D103").Value = 45072
$ws.Range("E103").Value = 15
$ws.Range("F103").Value = "Fruta"
$ws.Range("G103").Value = 100108
$ws.Range("H103").Value = "Tropicales y subtropicales"
$ws.Range("I103").Value = 100108003
$ws.Range("J103").Value = "Maracuyá"
$ws.Range("K103").Value = "Sin especificar"
$ws.Range("L103").Value = "Segunda"
$ws.Range("M103").Value = 160
$ws.Range("N103").Value = 23000
$ws.Range("O103").Value = 25000
$ws.Range("P103").Value = 24250
$ws.Range("Q103").Value = "$/caja 20 kilos"
$ws.Range("R103").Value = "Región de Arica y Parinacota"
$ws.Range("S103").Value = 1212
$ws.Range("T103").Value = 20
